# Update timestamps for the newly generated Handback report.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 161bd7f2-... on the Overview sheet
# (this value is also shown as de-de!H3, "Correspond Handoff Datetime").
$overview.Range("G3").Value = "2016-08-28 12:45:41"
$dede.Range("H3").Value = "2016-08-28 12:45:41"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 161bd7f2-...ab0372d42d13dd0fbd6ff201feb184b4d620e572.zh-cn.xlf
$zhcn.Range("H3").Value = "2016-08-28 12:45:37"
$zhcn.Range("K3").Value = "2016-08-28 12:46:22"

# de-de sheet: Correspond Handback DateTime
# for 161bd7f2-...ab0372d42d13dd0fbd6ff201feb184b4d620e572.de-de.xlf
$dede.Range("K3").Value = "2016-08-28 12:46:29"
